$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value = 'Golang Architect with API & Microservices'
$ws.Range("B52").Value = 'https://www.dice.com/job-detail/48983e1e-4ed5-423f-a666-706545588ed4'
$ws.Range("C52").Value = 'Remote'
$ws.Range("D52").Value = 'Contract, Third Party'
$ws.Range("E52").Value = '$75 - $80'
$ws.Range("F52").Value = 'Concent Software Solution LLC'

$ws.Range("A53").Value = 'Software Engineer ( Golang )'
$ws.Range("B53").Value = 'https://www.dice.com/job-detail/73d8e799-7f1d-4014-8e99-f8f0d221108a'
$ws.Range("C53").Value = 'Phoenix, Arizona'
$ws.Range("D53").Value = 'Contract'
$ws.Range("E53").Value = '$80 - $85'
$ws.Range("F53").Value = 'Source Mantra Inc'
